# Add 2022-Q4 data
#
# 1) Update the "总计" (totals) sheet: insert a new first data row for
#    2022-Q4 and shift the existing 2022-Q3 / 2022-Q2 rows down by one.
# 2) Insert a brand-new worksheet named "2022-Q4" right before the
#    existing "2022-Q3" worksheet, and fill it with the per-fund detail
#    rows for the new quarter (same shape as the existing quarter sheets).

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is currently active so that inserting the new
# worksheet below doesn't change which tab ends up marked as selected.
$activeSheetName = $wb.ActiveSheet.Name

# ---------------------------------------------------------------------
# 1. "总计" sheet (always the first sheet in the workbook)
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

# Push the current row 3 (2022-Q2) down into row 4, copying format first.
$totals.Cells.Item(3,1).Copy()
$totals.Cells.Item(4,1).PasteSpecial(-4122)   # xlPasteFormats

$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(4,2).Value = "2022-Q2"
$totals.Cells.Item(4,3).Value = 6
$totals.Cells.Item(4,4).Value = 1.29

# Current row 2 (2022-Q3) moves to row 3.
$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(3,2).Value = "2022-Q3"
$totals.Cells.Item(3,3).Value = 7
$totals.Cells.Item(3,4).Value = 0.82

# New 2022-Q4 row becomes row 2.
$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q4"
$totals.Cells.Item(2,3).Value = 10
$totals.Cells.Item(2,4).Value = 1.34

# ---------------------------------------------------------------------
# 2. New "2022-Q4" detail sheet, inserted before "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

for ($c = 2; $c -le 8; $c++) {
    $cell = $q4.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlVAlignTop
    $cell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cell.Borders.Item(7).Weight = 2
    $cell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
    $cell.Borders.Item(8).Weight = 2
    $cell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
    $cell.Borders.Item(9).Weight = 2
    $cell.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cell.Borders.Item(10).Weight = 2
}

# Columns B (fund code) and D:G (scale/position figures) hold their
# original text values verbatim (fund codes may have leading zeros, and
# the percentage-like figures are kept as text in the source data), so
# force a text number format before writing them.
$q4.Range("B2:B11").NumberFormat = "@"
$q4.Range("D2:G11").NumberFormat = "@"

$rows = @(
    @{ code = "202027"; name = "南方高端装备灵活配置混合A";   scale = "16.75"; pos = "93.09"; pct = "4.27"; mv = "0.7152"; rank = 8  },
    @{ code = "005207"; name = "南方高端装备灵活配置混合C";   scale = "5.35";  pos = "93.09"; pct = "4.27"; mv = "0.2284"; rank = 8  },
    @{ code = "000717"; name = "融通转型三动力灵活配置混合A"; scale = "3.27";  pos = "93.31"; pct = "4.80"; mv = "0.1570"; rank = 7  },
    @{ code = "161605"; name = "融通蓝筹成长混合";             scale = "4.89";  pos = "75.45"; pct = "3.07"; mv = "0.1501"; rank = 9  },
    @{ code = "004265"; name = "金鹰民丰回报定期开放混合";     scale = "4.49";  pos = "29.44"; pct = "0.93"; mv = "0.0418"; rank = 2  },
    @{ code = "011351"; name = "金鹰年年邮益一年持有期混合A"; scale = "3.04";  pos = "39.17"; pct = "1.18"; mv = "0.0359"; rank = 3  },
    @{ code = "016013"; name = "南方碳中和股票A";               scale = "0.42";  pos = "84.75"; pct = "2.96"; mv = "0.0124"; rank = 10 },
    @{ code = "011352"; name = "金鹰年年邮益一年持有期混合C"; scale = "0.23";  pos = "39.17"; pct = "1.18"; mv = "0.0027"; rank = 3  },
    @{ code = "009828"; name = "融通转型三动力灵活配置混合C"; scale = "0.02";  pos = "93.31"; pct = "4.80"; mv = "0.0010"; rank = 7  },
    @{ code = "016014"; name = "南方碳中和股票C";               scale = "0.01";  pos = "84.75"; pct = "2.96"; mv = "0.0003"; rank = 10 }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $q4.Cells.Item($r, 1).Value = $i
    $q4.Cells.Item($r, 2).Value = $row.code
    $q4.Cells.Item($r, 3).Value = $row.name
    $q4.Cells.Item($r, 4).Value = $row.scale
    $q4.Cells.Item($r, 5).Value = $row.pos
    $q4.Cells.Item($r, 6).Value = $row.pct
    $q4.Cells.Item($r, 7).Value = $row.mv
    $q4.Cells.Item($r, 8).Value = $row.rank
}

# Restore the original active sheet/tab selection (adding the new sheet
# above switches the active sheet to it).
$wb.Worksheets.Item($activeSheetName).Activate()

Write-Output "2022-Q4 sheet and totals updated"
